$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A64").Value = "UMB9LR"
$ws.Range("B64").Value = "Film de fusor HP"
$ws.Range("C64").Value = "1300 1320 Pro 400, M102 M104 M106 M125 M126 M127 M128 M130 M201 M203 M225 M227 M304 M305 M329 M400 M401 M402 M425 M426 M1136 M1212 M1213 M1216 M1217 M1536 M2727 P1000 P1005 P1010 P1012 P1015 P1018 P1022 P1050 P1102 P1106 P1150 P1160 P1200 P1220 P1606 P2014 P2015 P2030 P2035 P2055 P3010 P3011 P3015 P3016 P3018 P3020 P3030 P3050 P3055"
$ws.Range("D64").Value = 15000
$ws.Range("E64").Value = 100000
$ws.Range("F64").Value = 49
$ws.Range("G64").Value = 64
$ws.Range("H64").Formula = "=(E64-D64)*G64"
$ws.Range("I64").Formula = "=D64*F64"
$ws.Range("J64").Value = 735000
